# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1511
    3  = 29
    4  = 975
    6  = 2372
    8  = 1454
    10 = 161
    11 = 53
    12 = 403
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
